$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark so that it wraps the title paragraph
#    ("Miamly"): it should start right before "Miamly" and end right
#    after that paragraph's mark. Word keeps bookmark names unique, so
#    re-adding a bookmark named "_GoBack" here automatically removes the
#    old "_GoBack" bookmark further down in the document (around
#    "Toutes les astuces pour moins bien manger").
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$goBackRange = $d.Range($titlePara.End - 1, $titlePara.End - 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ---------------------------------------------------------------------
# 2) Split the run "<w:tab/><w:t>Recette au hasard</w:t>" into two runs:
#    one holding only the <w:tab/>, and a new one holding the text. Find
#    the tab immediately followed by the text, then rewrite that exact
#    span as two explicit runs via raw OOXML so the <w:tab/> element
#    (not a literal tab character) is preserved.
# ---------------------------------------------------------------------
$finder = $d.Content
$finder.Find.Execute("`tRecette au hasard")
$tileRange = $d.Range($finder.Start, $finder.End)
$tileXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr-FR"/></w:rPr><w:t>Recette au hasard</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tileRange.InsertXML($tileXml)

# ---------------------------------------------------------------------
# 3) Flag the "Normal Table" style as a quick style (adds <w:qFormat/>
#    to its <w:style> definition in styles.xml).
# ---------------------------------------------------------------------
$d.Styles.Item(4).QuickStyle = $true

Write-Output "edit complete"
